$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Row 2 ("Team ID") / Column 2: add "PNT2022TMID30874" and move the
# document's "_GoBack" bookmark to sit right after the new text. ---
$teamCell = $t.Cell(2, 2)
$teamRange = $teamCell.Range
$teamRange.Collapse(1)                       # wdCollapseStart: empty paragraph insertion point
$d.Bookmarks.Add("_GoBack", $teamRange)       # bookmark at the (still empty) insertion point
$teamInsert = $t.Cell(2, 2).Range
$teamInsert.Collapse(1)                       # wdCollapseStart again
$teamInsert.InsertBefore("PNT2022TMID30874")  # text lands before the bookmark -> bookmark ends up after it

# --- Row 3 ("Project Name") / Column 2: add the project name text. ---
$t2 = $d.Tables.Item(1)
$projCell = $t2.Cell(3, 2)
$projCell.Range.Text = "Intelligent vehicle damage and cost estimator"
